$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    # Force text storage so numeric-looking strings (e.g. "0.999")
    # are not coerced to Number type by COM, matching the source
    # workbook which stores these as inline strings.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "67.762.72"
Set-TextCell "D3" "3.324.77"
$ws.Range("E3").Value = "  -1.49%  "
Set-TextCell "D4" "0.999"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCell "D5" "581.14"
$ws.Range("E5").Value = "  -2.15%  "
Set-TextCell "D6" "173.79"
$ws.Range("E6").Value = "  -7.02%  "
Set-TextCell "D7" "0.999"
$ws.Range("E7").Value = "  -0.03%  "
Set-TextCell "D8" "0.585"
$ws.Range("E8").Value = "  -2.76%  "
Set-TextCell "D9" "3.317.53"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("E10").Value = "  -4.51%  "
Set-TextCell "D11" "0.576"
$ws.Range("E11").Value = "  -2.55%  "
Set-TextCell "D12" "45.43"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("E13").Value = "  -3.69%  "
Set-TextCell "D14" "665.50"
$ws.Range("E14").Value = "  +3.51%  "
Set-TextCell "D15" "3.854.66"
$ws.Range("E15").Value = "  -1.58%  "
Set-TextCell "D16" "8.38"
$ws.Range("E16").Value = "  -2.88%  "
Set-TextCell "D17" "67.832.24"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  -1.02%  "
Set-TextCell "D19" "3.313.50"
$ws.Range("E19").Value = "  -1.84%  "
Set-TextCell "D20" "17.43"
$ws.Range("E20").Value = "  -3.71%  "
Set-TextCell "D21" "10.88"
$ws.Range("E21").Value = "  -2.22%  "
Set-TextCell "D22" "0.887"
$ws.Range("E22").Value = "  -2.59%  "
Set-TextCell "D23" "5.39"
$ws.Range("E23").Value = "  +5.58%  "
Set-TextCell "D24" "16.99"
$ws.Range("E24").Value = "  -5.36%  "
Set-TextCell "D25" "97.42"
$ws.Range("E25").Value = "  -2.54%  "
Set-TextCell "D26" "3.84"
$ws.Range("E26").Value = "  -5.07%  "
Set-TextCell "D27" "2.67"
$ws.Range("E27").Value = "  -7.01%  "
Set-TextCell "D28" "9.28"
$ws.Range("E28").Value = "  -4.76%  "
Set-TextCell "D29" "33.44"
$ws.Range("E29").Value = "  +1.43%  "
Set-TextCell "D30" "8.39"
$ws.Range("E30").Value = "  -3.73%  "
Set-TextCell "D31" "7.32"
$ws.Range("E31").Value = "  +5.70%  "
Set-TextCell "D32" "583.78"
$ws.Range("E32").Value = "  -4.73%  "
Set-TextCell "D33" "10.94"
$ws.Range("E33").Value = "  -1.65%  "
Set-TextCell "D34" "0.104"
$ws.Range("E34").Value = "  -2.38%  "
$ws.Range("E35").Value = "  -0.06%  "
Set-TextCell "D36" "3.706.97"
$ws.Range("E36").Value = "  -8.45%  "
Set-TextCell "D37" "56.78"
$ws.Range("E37").Value = "  +0.73%  "
Set-TextCell "D38" "3.26"
$ws.Range("E38").Value = "  -14.19%  "
$ws.Range("E39").Value = "  +0.32%  "
Set-TextCell "D40" "32.62"
$ws.Range("E40").Value = "  -3.75%  "
Set-TextCell "D41" "2.61"
$ws.Range("E41").Value = "  -6.43%  "
Set-TextCell "D42" "3.09"
$ws.Range("E42").Value = "  -4.65%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D43" "0.332"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D44" "0.0₃0662"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell "D45" "3.24"
$ws.Range("E45").Value = "  -4.94%  "
Set-TextCell "D46" "0.0407"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("E50").Value = "  -4.06%  "
Set-TextCell "D51" "127.18"
$ws.Range("E51").Value = "  -0.70%  "
